$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap "Periodo Mora" text values for the two worker blocks (text-formatted cells)
$ws.Range("E16").Value = "2212"
$ws.Range("E17").Value = "2211"
$ws.Range("E18").Value = "2310"
$ws.Range("E19").Value = "2309"

# Swap "Valor Mora" amounts to match the swapped periods
$ws.Range("F18").Value = 46400
$ws.Range("F19").Value = 37120
